$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.717219114303589
$ws.Range("B1").Value = 2.257173299789429
$ws.Range("C1").Value = 3.34409499168396
$ws.Range("D1").Value = 3.881433963775635
$ws.Range("E1").Value = 0.6720036864280701
